$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits at the end of
#    the document (right after the "Scheduling" paragraph). Word moves this
#    special bookmark to the location of the most recent edit, so it needs to
#    be re-created at the new edit location (inside the title) below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Correct the title: "Tune Squad Meeting Minutes " becomes
#    "Tune Squad Teacher Meeting Minutes " - i.e. insert the word
#    "Teacher " between "Tune Squad " and "Meeting Minutes ".
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Meeting Minutes", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$insertionPoint = $titleRange.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore("Teacher ")

# ---------------------------------------------------------------------------
# 3. Force the title text into three distinct runs ("Tune Squad ", "Teacher ",
#    "Meeting Minutes ") by briefly dropping a bookmark at each split point.
#    The first split (between "Tune Squad " and "Teacher ") only needs a
#    transient bookmark; the second split is where the "_GoBack" bookmark
#    belongs permanently, matching the edit's insertion point.
# ---------------------------------------------------------------------------
$splitRange1 = $d.Content
$splitRange1.Find.Execute("Tune Squad ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$splitPoint1 = $splitRange1.Duplicate
$splitPoint1.Collapse(0)
$d.Bookmarks.Add("TempSplitMark", $splitPoint1)
$d.Bookmarks("TempSplitMark").Delete()

$splitRange2 = $d.Content
$splitRange2.Find.Execute("Teacher ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$splitPoint2 = $splitRange2.Duplicate
$splitPoint2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitPoint2)
